$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between C/E columns on the three "top face" rows (2, 8, 14)
$ws.Range("C2").Value = 3
$ws.Range("E2").Value = 1

$ws.Range("C8").Value = 3
$ws.Range("E8").Value = 1

$ws.Range("C14").Value = 8
$ws.Range("E14").Value = 3

# Update the active selection to match the author's final cursor position
$ws.Range("O15").Select()
